# Applies the "Add files via upload" edit to ga_setting.xlsx:
#  - parameter sheet:  A2:A6 renumbered 0..4 -> 1..5 (1-based index column)
#  - component sheet:  A2:A63 renumbered 0..61 -> 1..62 (1-based index column)
#  - setting sheet:    A1 gets the new label "Generative algorithm",
#                       D3 100 -> 30, D4 100 -> 15
#  - selections / active sheet restored to match the saved view state

$wb = $excel.ActiveWorkbook

$wsParameter = $wb.Worksheets.Item("parameter")
$wsComponent = $wb.Worksheets.Item("component")
$wsSetting   = $wb.Worksheets.Item("setting")

# --- parameter sheet: shift the index column A2:A6 up by one (0-based -> 1-based) ---
$wsParameter.Activate()
for ($r = 2; $r -le 6; $r++) {
    $wsParameter.Cells.Item($r, 1).Value = $r - 1
}
$wsParameter.Range("A2:A5").Select()

# --- component sheet: shift the index column A2:A63 up by one (0-based -> 1-based) ---
$wsComponent.Activate()
for ($r = 2; $r -le 63; $r++) {
    $wsComponent.Cells.Item($r, 1).Value = $r - 1
}
$wsComponent.Range("A52").Select()
$wsComponent.Range("A6").Select()

# --- setting sheet: new header label + updated generation counts ---
$wsSetting.Activate()
$wsSetting.Range("A1").Value = "Generative algorithm"
$wsSetting.Range("D3").Value = 30
$wsSetting.Range("D4").Value = 15
$wsSetting.Range("F8").Select()
